$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.082409381866455
$ws.Range("B1").Value = 2.731302261352539
$ws.Range("C1").Value = 1.967810273170471
$ws.Range("D1").Value = 1.814416646957397
$ws.Range("E1").Value = 1.869075655937195
